$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column "想去人数" values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 404
$ws1.Range("F4").Value = 451
$ws1.Range("F8").Value = 14203
$ws1.Range("F9").Value = 139
$ws1.Range("F10").Value = 106
$ws1.Range("F11").Value = 5705
$ws1.Range("F12").Value = 584
$ws1.Range("F19").Value = 174
$ws1.Range("F20").Value = 777
$ws1.Range("F21").Value = 2918
$ws1.Range("F23").Value = 10485
$ws1.Range("F25").Value = 48
$ws1.Range("F26").Value = 71
$ws1.Range("F27").Value = 3723
$ws1.Range("F28").Value = 241

# Sheet "全部类型" (sheet4): same events, shifted rows due to an extra
# "演出" entry present only in this aggregated sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 404
$ws4.Range("F5").Value = 451
$ws4.Range("F9").Value = 14203
$ws4.Range("F10").Value = 139
$ws4.Range("F11").Value = 106
$ws4.Range("F12").Value = 5705
$ws4.Range("F13").Value = 584
$ws4.Range("F20").Value = 174
$ws4.Range("F21").Value = 777
$ws4.Range("F22").Value = 2918
$ws4.Range("F25").Value = 10485
$ws4.Range("F27").Value = 48
$ws4.Range("F28").Value = 71
$ws4.Range("F29").Value = 3723
$ws4.Range("F30").Value = 241
